$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12236
$ws1.Range("F4").Value = 4496
$ws1.Range("F5").Value = 48
$ws1.Range("F6").Value = 67
$ws1.Range("F8").Value = 28
$ws1.Range("F9").Value = 2612
$ws1.Range("F11").Value = 206
$ws1.Range("F12").Value = 78
$ws1.Range("F13").Value = 5327
$ws1.Range("F15").Value = 213
$ws1.Range("F16").Value = 560
$ws1.Range("F17").Value = 11487
$ws1.Range("F18").Value = 11584
$ws1.Range("F20").Value = 65
$ws1.Range("F22").Value = 17

# Sheet "全部类型" (sheet4) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 12236
$ws4.Range("F4").Value = 4496
$ws4.Range("F5").Value = 48
$ws4.Range("F6").Value = 67
$ws4.Range("F8").Value = 28
$ws4.Range("F9").Value = 2612
$ws4.Range("F12").Value = 206
$ws4.Range("F13").Value = 78
$ws4.Range("F14").Value = 5327
$ws4.Range("F16").Value = 213
$ws4.Range("F17").Value = 560
$ws4.Range("F18").Value = 11487
$ws4.Range("F19").Value = 11584
$ws4.Range("F21").Value = 65
$ws4.Range("F23").Value = 17

